# ---------------------------------------------------------------------------
# Adapt column header formatting to respective input file names (#7)
#   * rename "<Column>_old"  -> "<Column>_FV2410"
#   * rename "<Column>_new"  -> "<Column>_FV2504"
#   * turn the data range into a proper Excel Table ("Table1")
#   * freeze the header row
# ---------------------------------------------------------------------------

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. rename the header cells (drives the shared-strings rename) --------
$fv2410 = @(
    "Segmentname_FV2410",
    "Segmentgruppe_FV2410",
    "Segment_FV2410",
    "Datenelement_FV2410",
    "Segment ID_FV2410",
    "Code_FV2410",
    "Qualifier_FV2410",
    "Beschreibung_FV2410",
    "Bedingungsausdruck_FV2410",
    "Bedingung_FV2410"
)
$fv2504 = @(
    "Segmentname_FV2504",
    "Segmentgruppe_FV2504",
    "Segment_FV2504",
    "Datenelement_FV2504",
    "Segment ID_FV2504",
    "Code_FV2504",
    "Qualifier_FV2504",
    "Beschreibung_FV2504",
    "Bedingungsausdruck_FV2504",
    "Bedingung_FV2504"
)

for ($i = 0; $i -lt $fv2410.Length; $i++) {
    $ws.Cells.Item(1, $i + 1).Value2 = $fv2410[$i]
}
# column 11 is "diff" - untouched
for ($i = 0; $i -lt $fv2504.Length; $i++) {
    $ws.Cells.Item(1, $i + 12).Value2 = $fv2504[$i]
}

# --- 2. turn A1:U76 into an Excel Table without Excel inventing a bold ----
#        header dxf (it only does that when the header range already
#        carries explicit formatting, which it does here) -----------------
$dataRange   = $ws.Range("A1:U76")
$headerRange = $ws.Range("A1:U1")
$stashRange  = $ws.Range("AA1:AU1")

$headerRange.Copy() | Out-Null
$stashRange.PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats) | Out-Null
$headerRange.ClearFormats() | Out-Null

$lo = $ws.ListObjects.Add(
    [Microsoft.Office.Interop.Excel.XlListObjectSourceType]::xlSrcRange,
    $dataRange,
    $null,
    [Microsoft.Office.Interop.Excel.XlYesNoGuess]::xlYes
)
$lo.Name = "Table1"
$lo.TableStyle = ""

# restore original header look now that the table exists
$stashRange.Copy() | Out-Null
$headerRange.PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats) | Out-Null
$stashRange.Clear() | Out-Null

# --- 3. freeze the header row ----------------------------------------------
$ws.Range("A2").Select() | Out-Null
$excel.ActiveWindow.FreezePanes = $true

Write-Host "Edit complete"
